# Applies a cyclic permutation of the observation records stored in rows
# 28, 29, 30, 32 and 33 of the "Artfynd" sheet.
#
#   row 28  <->  row 30            (2-cycle)
#   row 29  <-   row 32  <-  row 33 <- row 29   (3-cycle)
#
# Only the columns that actually carry per-record data move:
#   A (Id), B (Taxonsorteringsordning), D (Rödlistade), E (TaxonId),
#   F (Artnamn), G (Vetenskapligt namn), H (Auktor), Q (Ost), R (Nord).
# Columns C and I are identical across these rows, so they are left as-is.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A", "B", "D", "E", "F", "G", "H", "Q", "R")

function Get-RowData($ws, $row, $cols) {
    $data = @{}
    foreach ($c in $cols) {
        $data[$c] = $ws.Range("$c$row").Value2
    }
    return $data
}

function Set-RowData($ws, $row, $cols, $data) {
    foreach ($c in $cols) {
        $ws.Range("$c$row").Value = $data[$c]
    }
}

# Snapshot the current ("before") state of every row involved.
$row28 = Get-RowData $ws 28 $cols
$row29 = Get-RowData $ws 29 $cols
$row30 = Get-RowData $ws 30 $cols
$row32 = Get-RowData $ws 32 $cols
$row33 = Get-RowData $ws 33 $cols

# Write back the permuted ("after") state.
Set-RowData $ws 28 $cols $row30
Set-RowData $ws 29 $cols $row32
Set-RowData $ws 30 $cols $row28
Set-RowData $ws 32 $cols $row33
Set-RowData $ws 33 $cols $row29
